# GossA-HW40: re-ran simulation.
#   - two brand new structures ("Holden", "Rizzie Spiral") were inserted
#     right after "Spiral5", pushing every later category down by two rows
#   - "Thomas Hex" was renamed to "Matthies Hex"
#   - the whole table grew from 29 to 31 data rows (A1:W29 -> A1:W31)
#
# Implementation notes:
#   Column A is always just the zero-based row index (row-2); it does not
#   "shift" - it is simply re-sequenced after growing the table.
#   Columns B:W (category label + simulation numbers) for the rows that
#   already existed all move down by two rows, unchanged. Rows 4 and 5 are
#   populated with freshly computed simulation numbers for the two new
#   categories.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Capture the existing label+data block (old rows 4..29, columns B:W)
#    BEFORE anything is overwritten, then shift it down two rows so it
#    lands at rows 6..31.
# ---------------------------------------------------------------------
$srcBlock = $ws.Range("B4:W29").Value()
$ws.Range("B6:W31").Value = $srcBlock

# ---------------------------------------------------------------------
# 2) Re-sequence column A (0-based row index) for every data row,
#    including the two brand new rows at the bottom (30 and 31).
# ---------------------------------------------------------------------
for ($r = 2; $r -le 31; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 2
}

# Column A carries bold+border formatting (same style as row 2..29's A
# cells); copy that formatting down onto the two newly-created rows so we
# don't end up with an unformatted cell there.
$ws.Range("A29").Copy()
$ws.Range("A30:A31").PasteSpecial(-4122)  # xlPasteFormats
$ws.Cells.Item(30, 1).Value = 28
$ws.Cells.Item(31, 1).Value = 29

# ---------------------------------------------------------------------
# 3) Write the two new categories into rows 4 and 5 with their freshly
#    computed simulation results.
# ---------------------------------------------------------------------
$cols = @("C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W")

$ws.Range("B4").Value = "Holden"
$row4 = @(1.00498667424703,1.00498667424703,0.9983515112251989,1.002864485424185,1.001878931611703,0.9950110452929396,0.9923057148114121,0.9985125696001621,1.002003925322815,0.9923057148114121,1.00498667424703,1.00498667424703,1.002003925322815,0.9971548200671135,1.000177718274007,0.9997654381270856,0.9975537171198087,0.9997654381270856,0.999411956401614,1.000526899970697,0.9994893571919307)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Cells.Item(4, 3 + $i).Value = $row4[$i]
}

$ws.Range("B5").Value = "Rizzie Spiral"
$row5 = @(1.010338528141337,1.010338528141337,0.9968226724616092,1.005425666555968,1.00384338652565,0.9903881835754523,0.9852977042605678,0.9968662305875478,1.00337351694334,0.9852977042605678,1.010338528141337,1.010338528141337,1.00337351694334,0.9943356106019539,1.000098094702474,0.9996699164484149,0.9951646312218388,0.9996699164484149,0.9989581054517134,1.001234189989638,0.9990444861314339)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Cells.Item(5, 3 + $i).Value = $row5[$i]
}

# ---------------------------------------------------------------------
# 4) Rename "Thomas Hex" -> "Matthies Hex" (now living a couple of rows
#    further down the table thanks to the insert above).
# ---------------------------------------------------------------------
for ($r = 6; $r -le 31; $r++) {
    $cell = $ws.Cells.Item($r, 2)
    if ($cell.Value() -eq "Thomas Hex") {
        $cell.Value = "Matthies Hex"
    }
}
